$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns before the existing "general_college_subjects.arts"
# column (currently column R). This shifts R:AE -> U:AH.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the inserted columns.
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for row 2 in the inserted columns.
$ws.Range("R2").Value = 3
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 1

# Normalize casing of the descriptive importance/consideration values in row 2.
$ws.Range("D2").Value = "not considered"
$ws.Range("E2").Value = "important"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
